$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill constant columns (A,B,C,E,F,G,H,R) for new rows 165-167
$ws.Range("A165").Value = 9
$ws.Range("B165").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C165").Value = "Metropolitana"
$ws.Range("E165").Value = 13
$ws.Range("F165").Value = 100112052
$ws.Range("G165").Value = "Albahaca"
$ws.Range("H165").Value = "Sin especificar"
$ws.Range("R165").Value = "Hortaliza"
$ws.Range("D165").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A166").Value = 9
$ws.Range("B166").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C166").Value = "Metropolitana"
$ws.Range("E166").Value = 13
$ws.Range("F166").Value = 100112052
$ws.Range("G166").Value = "Albahaca"
$ws.Range("H166").Value = "Sin especificar"
$ws.Range("R166").Value = "Hortaliza"
$ws.Range("D166").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A167").Value = 9
$ws.Range("B167").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C167").Value = "Metropolitana"
$ws.Range("E167").Value = 13
$ws.Range("F167").Value = 100112052
$ws.Range("G167").Value = "Albahaca"
$ws.Range("H167").Value = "Sin especificar"
$ws.Range("R167").Value = "Hortaliza"
$ws.Range("D167").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Update/fill variable columns D,I,J,K,L,M,N,O,P,Q for rows 69-167
$ws.Range("D69").Value = (Get-Date -Year 2021 -Month 8 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 250
$ws.Range("K69").Value = 5000
$ws.Range("L69").Value = 6000
$ws.Range("M69").Value = 5500
$ws.Range("N69").Value = "`$/paquete"
$ws.Range("O69").Value = "Región de Arica y Parinacota"
$ws.Range("P69").Value = 5500
$ws.Range("Q69").Value = 1
$ws.Range("D70").Value = (Get-Date -Year 2021 -Month 8 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I70").Value = "Segunda"
$ws.Range("J70").Value = 160
$ws.Range("K70").Value = 4500
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = 4500
$ws.Range("N70").Value = "`$/paquete"
$ws.Range("O70").Value = "Región de Arica y Parinacota"
$ws.Range("P70").Value = 4500
$ws.Range("Q70").Value = 1
$ws.Range("D71").Value = (Get-Date -Year 2020 -Month 12 -Day 11 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I71").Value = "Primera"
$ws.Range("J71").Value = 160
$ws.Range("K71").Value = 4500
$ws.Range("L71").Value = 5000
$ws.Range("M71").Value = 4750
$ws.Range("N71").Value = "`$/docena de matas"
$ws.Range("O71").Value = "Región Metropolitana"
$ws.Range("P71").Value = 792
$ws.Range("Q71").Value = 6
$ws.Range("D72").Value = (Get-Date -Year 2021 -Month 1 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 490
$ws.Range("K72").Value = 3000
$ws.Range("L72").Value = 4000
$ws.Range("M72").Value = 3490
$ws.Range("N72").Value = "`$/docena de matas"
$ws.Range("O72").Value = "Región Metropolitana"
$ws.Range("P72").Value = 582
$ws.Range("Q72").Value = 6
$ws.Range("D73").Value = (Get-Date -Year 2021 -Month 2 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I73").Value = "Primera"
$ws.Range("J73").Value = 450
$ws.Range("K73").Value = 2500
$ws.Range("L73").Value = 3000
$ws.Range("M73").Value = 2778
$ws.Range("N73").Value = "`$/docena de matas"
$ws.Range("O73").Value = "Provincia de Chacabuco"
$ws.Range("P73").Value = 463
$ws.Range("Q73").Value = 6
$ws.Range("D74").Value = (Get-Date -Year 2021 -Month 3 -Day 18 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I74").Value = "Primera"
$ws.Range("J74").Value = 340
$ws.Range("K74").Value = 3500
$ws.Range("L74").Value = 3500
$ws.Range("M74").Value = 3500
$ws.Range("N74").Value = "`$/docena de matas"
$ws.Range("O74").Value = "Región Metropolitana"
$ws.Range("P74").Value = 583
$ws.Range("Q74").Value = 6
$ws.Range("D75").Value = (Get-Date -Year 2021 -Month 3 -Day 18 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I75").Value = "Segunda"
$ws.Range("J75").Value = 250
$ws.Range("K75").Value = 2500
$ws.Range("L75").Value = 2500
$ws.Range("M75").Value = 2500
$ws.Range("N75").Value = "`$/docena de matas"
$ws.Range("O75").Value = "Región Metropolitana"
$ws.Range("P75").Value = 417
$ws.Range("Q75").Value = 6
$ws.Range("D76").Value = (Get-Date -Year 2021 -Month 3 -Day 2 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I76").Value = "Primera"
$ws.Range("J76").Value = 340
$ws.Range("K76").Value = 3000
$ws.Range("L76").Value = 3000
$ws.Range("M76").Value = 3000
$ws.Range("N76").Value = "`$/docena de matas"
$ws.Range("O76").Value = "Región Metropolitana"
$ws.Range("P76").Value = 500
$ws.Range("Q76").Value = 6
$ws.Range("D77").Value = (Get-Date -Year 2021 -Month 1 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I77").Value = "Primera"
$ws.Range("J77").Value = 480
$ws.Range("K77").Value = 3000
$ws.Range("L77").Value = 3000
$ws.Range("M77").Value = 3000
$ws.Range("N77").Value = "`$/docena de matas"
$ws.Range("O77").Value = "Región Metropolitana"
$ws.Range("P77").Value = 500
$ws.Range("Q77").Value = 6
$ws.Range("D78").Value = (Get-Date -Year 2021 -Month 1 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 190
$ws.Range("K78").Value = 3000
$ws.Range("L78").Value = 3000
$ws.Range("M78").Value = 3000
$ws.Range("N78").Value = "`$/docena de matas"
$ws.Range("O78").Value = "Región de O'Higgins"
$ws.Range("P78").Value = 500
$ws.Range("Q78").Value = 6
$ws.Range("D79").Value = (Get-Date -Year 2021 -Month 4 -Day 9 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 340
$ws.Range("K79").Value = 3500
$ws.Range("L79").Value = 3500
$ws.Range("M79").Value = 3500
$ws.Range("N79").Value = "`$/docena de matas"
$ws.Range("O79").Value = "Región Metropolitana"
$ws.Range("P79").Value = 583
$ws.Range("Q79").Value = 6
$ws.Range("D80").Value = (Get-Date -Year 2021 -Month 4 -Day 9 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I80").Value = "Segunda"
$ws.Range("J80").Value = 250
$ws.Range("K80").Value = 2500
$ws.Range("L80").Value = 2500
$ws.Range("M80").Value = 2500
$ws.Range("N80").Value = "`$/docena de matas"
$ws.Range("O80").Value = "Región Metropolitana"
$ws.Range("P80").Value = 417
$ws.Range("Q80").Value = 6
$ws.Range("D81").Value = (Get-Date -Year 2020 -Month 12 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I81").Value = "Primera"
$ws.Range("J81").Value = 430
$ws.Range("K81").Value = 3500
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = 3750
$ws.Range("N81").Value = "`$/docena de matas"
$ws.Range("O81").Value = "Región Metropolitana"
$ws.Range("P81").Value = 625
$ws.Range("Q81").Value = 6
$ws.Range("D82").Value = (Get-Date -Year 2020 -Month 12 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I82").Value = "Segunda"
$ws.Range("J82").Value = 250
$ws.Range("K82").Value = 3000
$ws.Range("L82").Value = 3000
$ws.Range("M82").Value = 3000
$ws.Range("N82").Value = "`$/docena de matas"
$ws.Range("O82").Value = "Región Metropolitana"
$ws.Range("P82").Value = 500
$ws.Range("Q82").Value = 6
$ws.Range("D83").Value = (Get-Date -Year 2021 -Month 6 -Day 7 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 260
$ws.Range("K83").Value = 4500
$ws.Range("L83").Value = 5000
$ws.Range("M83").Value = 4750
$ws.Range("N83").Value = "`$/paquete"
$ws.Range("O83").Value = "Región de Arica y Parinacota"
$ws.Range("P83").Value = 4750
$ws.Range("Q83").Value = 1
$ws.Range("D84").Value = (Get-Date -Year 2020 -Month 12 -Day 4 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 70
$ws.Range("K84").Value = 5000
$ws.Range("L84").Value = 5000
$ws.Range("M84").Value = 5000
$ws.Range("N84").Value = "`$/docena de matas"
$ws.Range("O84").Value = "Provincia de Quillota"
$ws.Range("P84").Value = 833
$ws.Range("Q84").Value = 6
$ws.Range("D85").Value = (Get-Date -Year 2020 -Month 12 -Day 4 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I85").Value = "Primera"
$ws.Range("J85").Value = 250
$ws.Range("K85").Value = 4500
$ws.Range("L85").Value = 5000
$ws.Range("M85").Value = 4800
$ws.Range("N85").Value = "`$/docena de matas"
$ws.Range("O85").Value = "Región Metropolitana"
$ws.Range("P85").Value = 800
$ws.Range("Q85").Value = 6
$ws.Range("D86").Value = (Get-Date -Year 2021 -Month 8 -Day 2 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I86").Value = "Primera"
$ws.Range("J86").Value = 193
$ws.Range("K86").Value = 5000
$ws.Range("L86").Value = 6000
$ws.Range("M86").Value = 5497
$ws.Range("N86").Value = "`$/paquete"
$ws.Range("O86").Value = "Región de Arica y Parinacota"
$ws.Range("P86").Value = 5497
$ws.Range("Q86").Value = 1
$ws.Range("D87").Value = (Get-Date -Year 2021 -Month 5 -Day 20 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I87").Value = "Primera"
$ws.Range("J87").Value = 250
$ws.Range("K87").Value = 3500
$ws.Range("L87").Value = 4000
$ws.Range("M87").Value = 3750
$ws.Range("N87").Value = "`$/paquete"
$ws.Range("O87").Value = "Región de Arica y Parinacota"
$ws.Range("P87").Value = 3750
$ws.Range("Q87").Value = 1
$ws.Range("D88").Value = (Get-Date -Year 2021 -Month 4 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 340
$ws.Range("K88").Value = 3500
$ws.Range("L88").Value = 3500
$ws.Range("M88").Value = 3500
$ws.Range("N88").Value = "`$/docena de matas"
$ws.Range("O88").Value = "Región Metropolitana"
$ws.Range("P88").Value = 583
$ws.Range("Q88").Value = 6
$ws.Range("D89").Value = (Get-Date -Year 2021 -Month 4 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I89").Value = "Segunda"
$ws.Range("J89").Value = 160
$ws.Range("K89").Value = 2500
$ws.Range("L89").Value = 2500
$ws.Range("M89").Value = 2500
$ws.Range("N89").Value = "`$/docena de matas"
$ws.Range("O89").Value = "Región Metropolitana"
$ws.Range("P89").Value = 417
$ws.Range("Q89").Value = 6
$ws.Range("D90").Value = (Get-Date -Year 2020 -Month 12 -Day 7 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I90").Value = "Primera"
$ws.Range("J90").Value = 430
$ws.Range("K90").Value = 3000
$ws.Range("L90").Value = 3500
$ws.Range("M90").Value = 3250
$ws.Range("N90").Value = "`$/paquete"
$ws.Range("O90").Value = "Perú"
$ws.Range("P90").Value = 3250
$ws.Range("Q90").Value = 1
$ws.Range("D91").Value = (Get-Date -Year 2020 -Month 12 -Day 7 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I91").Value = "Primera"
$ws.Range("J91").Value = 340
$ws.Range("K91").Value = 2500
$ws.Range("L91").Value = 3000
$ws.Range("M91").Value = 2750
$ws.Range("N91").Value = "`$/paquete"
$ws.Range("O91").Value = "Región de Arica y Parinacota"
$ws.Range("P91").Value = 2750
$ws.Range("Q91").Value = 1
$ws.Range("D92").Value = (Get-Date -Year 2021 -Month 1 -Day 18 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I92").Value = "Primera"
$ws.Range("J92").Value = 150
$ws.Range("K92").Value = 4000
$ws.Range("L92").Value = 4000
$ws.Range("M92").Value = 4000
$ws.Range("N92").Value = "`$/docena de matas"
$ws.Range("O92").Value = "Región Metropolitana"
$ws.Range("P92").Value = 667
$ws.Range("Q92").Value = 6
$ws.Range("D93").Value = (Get-Date -Year 2021 -Month 6 -Day 24 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I93").Value = "Primera"
$ws.Range("J93").Value = 270
$ws.Range("K93").Value = 4500
$ws.Range("L93").Value = 5000
$ws.Range("M93").Value = 4750
$ws.Range("N93").Value = "`$/paquete"
$ws.Range("O93").Value = "Región de Arica y Parinacota"
$ws.Range("P93").Value = 4750
$ws.Range("Q93").Value = 1
$ws.Range("D94").Value = (Get-Date -Year 2021 -Month 6 -Day 24 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I94").Value = "Segunda"
$ws.Range("J94").Value = 160
$ws.Range("K94").Value = 3500
$ws.Range("L94").Value = 3500
$ws.Range("M94").Value = 3500
$ws.Range("N94").Value = "`$/paquete"
$ws.Range("O94").Value = "Región de Arica y Parinacota"
$ws.Range("P94").Value = 3500
$ws.Range("Q94").Value = 1
$ws.Range("D95").Value = (Get-Date -Year 2021 -Month 2 -Day 11 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I95").Value = "Primera"
$ws.Range("J95").Value = 700
$ws.Range("K95").Value = 2500
$ws.Range("L95").Value = 3000
$ws.Range("M95").Value = 2750
$ws.Range("N95").Value = "`$/docena de matas"
$ws.Range("O95").Value = "Provincia de Chacabuco"
$ws.Range("P95").Value = 458
$ws.Range("Q95").Value = 6
$ws.Range("D96").Value = (Get-Date -Year 2021 -Month 2 -Day 11 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I96").Value = "Primera"
$ws.Range("J96").Value = 160
$ws.Range("K96").Value = 3000
$ws.Range("L96").Value = 3000
$ws.Range("M96").Value = 3000
$ws.Range("N96").Value = "`$/docena de matas"
$ws.Range("O96").Value = "Región de O'Higgins"
$ws.Range("P96").Value = 500
$ws.Range("Q96").Value = 6
$ws.Range("D97").Value = (Get-Date -Year 2021 -Month 3 -Day 19 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I97").Value = "Primera"
$ws.Range("J97").Value = 340
$ws.Range("K97").Value = 3500
$ws.Range("L97").Value = 3500
$ws.Range("M97").Value = 3500
$ws.Range("N97").Value = "`$/docena de matas"
$ws.Range("O97").Value = "Región Metropolitana"
$ws.Range("P97").Value = 583
$ws.Range("Q97").Value = 6
$ws.Range("D98").Value = (Get-Date -Year 2021 -Month 3 -Day 19 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I98").Value = "Segunda"
$ws.Range("J98").Value = 250
$ws.Range("K98").Value = 2500
$ws.Range("L98").Value = 2500
$ws.Range("M98").Value = 2500
$ws.Range("N98").Value = "`$/docena de matas"
$ws.Range("O98").Value = "Región Metropolitana"
$ws.Range("P98").Value = 417
$ws.Range("Q98").Value = 6
$ws.Range("D99").Value = (Get-Date -Year 2021 -Month 4 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 250
$ws.Range("K99").Value = 5000
$ws.Range("L99").Value = 5000
$ws.Range("M99").Value = 5000
$ws.Range("N99").Value = "`$/docena de matas"
$ws.Range("O99").Value = "Región Metropolitana"
$ws.Range("P99").Value = 833
$ws.Range("Q99").Value = 6
$ws.Range("D100").Value = (Get-Date -Year 2021 -Month 6 -Day 17 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I100").Value = "Primera"
$ws.Range("J100").Value = 220
$ws.Range("K100").Value = 4500
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = 4750
$ws.Range("N100").Value = "`$/paquete"
$ws.Range("O100").Value = "Región de Arica y Parinacota"
$ws.Range("P100").Value = 4750
$ws.Range("Q100").Value = 1
$ws.Range("D101").Value = (Get-Date -Year 2021 -Month 1 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 380
$ws.Range("K101").Value = 3000
$ws.Range("L101").Value = 3000
$ws.Range("M101").Value = 3000
$ws.Range("N101").Value = "`$/docena de matas"
$ws.Range("O101").Value = "Provincia de Chacabuco"
$ws.Range("P101").Value = 500
$ws.Range("Q101").Value = 6
$ws.Range("D102").Value = (Get-Date -Year 2021 -Month 1 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I102").Value = "Primera"
$ws.Range("J102").Value = 120
$ws.Range("K102").Value = 3000
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = 3000
$ws.Range("N102").Value = "`$/docena de matas"
$ws.Range("O102").Value = "Región de O'Higgins"
$ws.Range("P102").Value = 500
$ws.Range("Q102").Value = 6
$ws.Range("D103").Value = (Get-Date -Year 2021 -Month 3 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I103").Value = "Primera"
$ws.Range("J103").Value = 340
$ws.Range("K103").Value = 3500
$ws.Range("L103").Value = 3500
$ws.Range("M103").Value = 3500
$ws.Range("N103").Value = "`$/docena de matas"
$ws.Range("O103").Value = "Región Metropolitana"
$ws.Range("P103").Value = 583
$ws.Range("Q103").Value = 6
$ws.Range("D104").Value = (Get-Date -Year 2021 -Month 1 -Day 13 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I104").Value = "Primera"
$ws.Range("J104").Value = 380
$ws.Range("K104").Value = 3000
$ws.Range("L104").Value = 4000
$ws.Range("M104").Value = 3526
$ws.Range("N104").Value = "`$/docena de matas"
$ws.Range("O104").Value = "Región Metropolitana"
$ws.Range("P104").Value = 588
$ws.Range("Q104").Value = 6
$ws.Range("D105").Value = (Get-Date -Year 2021 -Month 2 -Day 4 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I105").Value = "Primera"
$ws.Range("J105").Value = 400
$ws.Range("K105").Value = 3000
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = 3375
$ws.Range("N105").Value = "`$/docena de matas"
$ws.Range("O105").Value = "Provincia de Chacabuco"
$ws.Range("P105").Value = 562
$ws.Range("Q105").Value = 6
$ws.Range("D106").Value = (Get-Date -Year 2021 -Month 4 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I106").Value = "Primera"
$ws.Range("J106").Value = 160
$ws.Range("K106").Value = 5000
$ws.Range("L106").Value = 5000
$ws.Range("M106").Value = 5000
$ws.Range("N106").Value = "`$/docena de matas"
$ws.Range("O106").Value = "Región Metropolitana"
$ws.Range("P106").Value = 833
$ws.Range("Q106").Value = 6
$ws.Range("D107").Value = (Get-Date -Year 2020 -Month 12 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I107").Value = "Primera"
$ws.Range("J107").Value = 120
$ws.Range("K107").Value = 4000
$ws.Range("L107").Value = 4000
$ws.Range("M107").Value = 4000
$ws.Range("N107").Value = "`$/docena de matas"
$ws.Range("O107").Value = "Provincia de Chacabuco"
$ws.Range("P107").Value = 667
$ws.Range("Q107").Value = 6
$ws.Range("D108").Value = (Get-Date -Year 2021 -Month 6 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I108").Value = "Primera"
$ws.Range("J108").Value = 210
$ws.Range("K108").Value = 3500
$ws.Range("L108").Value = 4000
$ws.Range("M108").Value = 3750
$ws.Range("N108").Value = "`$/paquete"
$ws.Range("O108").Value = "Región de Arica y Parinacota"
$ws.Range("P108").Value = 3750
$ws.Range("Q108").Value = 1
$ws.Range("D109").Value = (Get-Date -Year 2020 -Month 12 -Day 31 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I109").Value = "Primera"
$ws.Range("J109").Value = 340
$ws.Range("K109").Value = 4500
$ws.Range("L109").Value = 5000
$ws.Range("M109").Value = 4750
$ws.Range("N109").Value = "`$/docena de matas"
$ws.Range("O109").Value = "Región Metropolitana"
$ws.Range("P109").Value = 792
$ws.Range("Q109").Value = 6
$ws.Range("D110").Value = (Get-Date -Year 2020 -Month 12 -Day 31 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I110").Value = "Segunda"
$ws.Range("J110").Value = 160
$ws.Range("K110").Value = 7000
$ws.Range("L110").Value = 7000
$ws.Range("M110").Value = 7000
$ws.Range("N110").Value = "`$/docena de matas"
$ws.Range("O110").Value = "Región Metropolitana"
$ws.Range("P110").Value = 1167
$ws.Range("Q110").Value = 6
$ws.Range("D111").Value = (Get-Date -Year 2021 -Month 4 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I111").Value = "Primera"
$ws.Range("J111").Value = 250
$ws.Range("K111").Value = 5500
$ws.Range("L111").Value = 5500
$ws.Range("M111").Value = 5500
$ws.Range("N111").Value = "`$/docena de matas"
$ws.Range("O111").Value = "Región Metropolitana"
$ws.Range("P111").Value = 917
$ws.Range("Q111").Value = 6
$ws.Range("D112").Value = (Get-Date -Year 2021 -Month 4 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I112").Value = "Segunda"
$ws.Range("J112").Value = 160
$ws.Range("K112").Value = 4500
$ws.Range("L112").Value = 4500
$ws.Range("M112").Value = 4500
$ws.Range("N112").Value = "`$/docena de matas"
$ws.Range("O112").Value = "Región Metropolitana"
$ws.Range("P112").Value = 750
$ws.Range("Q112").Value = 6
$ws.Range("D113").Value = (Get-Date -Year 2021 -Month 3 -Day 23 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I113").Value = "Primera"
$ws.Range("J113").Value = 340
$ws.Range("K113").Value = 3500
$ws.Range("L113").Value = 3500
$ws.Range("M113").Value = 3500
$ws.Range("N113").Value = "`$/docena de matas"
$ws.Range("O113").Value = "Región Metropolitana"
$ws.Range("P113").Value = 583
$ws.Range("Q113").Value = 6
$ws.Range("D114").Value = (Get-Date -Year 2021 -Month 3 -Day 23 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I114").Value = "Segunda"
$ws.Range("J114").Value = 250
$ws.Range("K114").Value = 2500
$ws.Range("L114").Value = 2500
$ws.Range("M114").Value = 2500
$ws.Range("N114").Value = "`$/docena de matas"
$ws.Range("O114").Value = "Región Metropolitana"
$ws.Range("P114").Value = 417
$ws.Range("Q114").Value = 6
$ws.Range("D115").Value = (Get-Date -Year 2021 -Month 7 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I115").Value = "Primera"
$ws.Range("J115").Value = 250
$ws.Range("K115").Value = 4500
$ws.Range("L115").Value = 5000
$ws.Range("M115").Value = 4750
$ws.Range("N115").Value = "`$/paquete"
$ws.Range("O115").Value = "Región de Arica y Parinacota"
$ws.Range("P115").Value = 4750
$ws.Range("Q115").Value = 1
$ws.Range("D116").Value = (Get-Date -Year 2021 -Month 7 -Day 22 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I116").Value = "Primera"
$ws.Range("J116").Value = 250
$ws.Range("K116").Value = 4500
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = 4750
$ws.Range("N116").Value = "`$/paquete"
$ws.Range("O116").Value = "Región de Arica y Parinacota"
$ws.Range("P116").Value = 4750
$ws.Range("Q116").Value = 1
$ws.Range("D117").Value = (Get-Date -Year 2021 -Month 7 -Day 22 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I117").Value = "Segunda"
$ws.Range("J117").Value = 160
$ws.Range("K117").Value = 4000
$ws.Range("L117").Value = 4000
$ws.Range("M117").Value = 4000
$ws.Range("N117").Value = "`$/paquete"
$ws.Range("O117").Value = "Región de Arica y Parinacota"
$ws.Range("P117").Value = 4000
$ws.Range("Q117").Value = 1
$ws.Range("D118").Value = (Get-Date -Year 2021 -Month 2 -Day 24 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I118").Value = "Primera"
$ws.Range("J118").Value = 340
$ws.Range("K118").Value = 3000
$ws.Range("L118").Value = 3500
$ws.Range("M118").Value = 3250
$ws.Range("N118").Value = "`$/docena de matas"
$ws.Range("O118").Value = "Región Metropolitana"
$ws.Range("P118").Value = 542
$ws.Range("Q118").Value = 6
$ws.Range("D119").Value = (Get-Date -Year 2021 -Month 2 -Day 23 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I119").Value = "Primera"
$ws.Range("J119").Value = 340
$ws.Range("K119").Value = 3000
$ws.Range("L119").Value = 3000
$ws.Range("M119").Value = 3000
$ws.Range("N119").Value = "`$/docena de matas"
$ws.Range("O119").Value = "Región Metropolitana"
$ws.Range("P119").Value = 500
$ws.Range("Q119").Value = 6
$ws.Range("D120").Value = (Get-Date -Year 2021 -Month 4 -Day 6 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I120").Value = "Primera"
$ws.Range("J120").Value = 340
$ws.Range("K120").Value = 3500
$ws.Range("L120").Value = 3500
$ws.Range("M120").Value = 3500
$ws.Range("N120").Value = "`$/docena de matas"
$ws.Range("O120").Value = "Región Metropolitana"
$ws.Range("P120").Value = 583
$ws.Range("Q120").Value = 6
$ws.Range("D121").Value = (Get-Date -Year 2021 -Month 4 -Day 6 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I121").Value = "Segunda"
$ws.Range("J121").Value = 160
$ws.Range("K121").Value = 2500
$ws.Range("L121").Value = 2500
$ws.Range("M121").Value = 2500
$ws.Range("N121").Value = "`$/docena de matas"
$ws.Range("O121").Value = "Región Metropolitana"
$ws.Range("P121").Value = 417
$ws.Range("Q121").Value = 6
$ws.Range("D122").Value = (Get-Date -Year 2021 -Month 4 -Day 8 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I122").Value = "Primera"
$ws.Range("J122").Value = 340
$ws.Range("K122").Value = 3500
$ws.Range("L122").Value = 3500
$ws.Range("M122").Value = 3500
$ws.Range("N122").Value = "`$/docena de matas"
$ws.Range("O122").Value = "Región Metropolitana"
$ws.Range("P122").Value = 583
$ws.Range("Q122").Value = 6
$ws.Range("D123").Value = (Get-Date -Year 2021 -Month 4 -Day 8 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I123").Value = "Segunda"
$ws.Range("J123").Value = 250
$ws.Range("K123").Value = 2500
$ws.Range("L123").Value = 2500
$ws.Range("M123").Value = 2500
$ws.Range("N123").Value = "`$/docena de matas"
$ws.Range("O123").Value = "Región Metropolitana"
$ws.Range("P123").Value = 417
$ws.Range("Q123").Value = 6
$ws.Range("D124").Value = (Get-Date -Year 2021 -Month 2 -Day 18 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I124").Value = "Primera"
$ws.Range("J124").Value = 430
$ws.Range("K124").Value = 3000
$ws.Range("L124").Value = 3000
$ws.Range("M124").Value = 3000
$ws.Range("N124").Value = "`$/docena de matas"
$ws.Range("O124").Value = "Región Metropolitana"
$ws.Range("P124").Value = 500
$ws.Range("Q124").Value = 6
$ws.Range("D125").Value = (Get-Date -Year 2021 -Month 1 -Day 6 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I125").Value = "Primera"
$ws.Range("J125").Value = 330
$ws.Range("K125").Value = 3500
$ws.Range("L125").Value = 4000
$ws.Range("M125").Value = 3803
$ws.Range("N125").Value = "`$/docena de matas"
$ws.Range("O125").Value = "Provincia de Chacabuco"
$ws.Range("P125").Value = 634
$ws.Range("Q125").Value = 6
$ws.Range("D126").Value = (Get-Date -Year 2020 -Month 12 -Day 16 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I126").Value = "Primera"
$ws.Range("J126").Value = 260
$ws.Range("K126").Value = 4000
$ws.Range("L126").Value = 5000
$ws.Range("M126").Value = 4385
$ws.Range("N126").Value = "`$/docena de matas"
$ws.Range("O126").Value = "Provincia de Chacabuco"
$ws.Range("P126").Value = 731
$ws.Range("Q126").Value = 6
$ws.Range("D127").Value = (Get-Date -Year 2021 -Month 6 -Day 22 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I127").Value = "Primera"
$ws.Range("J127").Value = 250
$ws.Range("K127").Value = 4500
$ws.Range("L127").Value = 5000
$ws.Range("M127").Value = 4750
$ws.Range("N127").Value = "`$/paquete"
$ws.Range("O127").Value = "Región de Arica y Parinacota"
$ws.Range("P127").Value = 4750
$ws.Range("Q127").Value = 1
$ws.Range("D128").Value = (Get-Date -Year 2020 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I128").Value = "Primera"
$ws.Range("J128").Value = 340
$ws.Range("K128").Value = 4500
$ws.Range("L128").Value = 5000
$ws.Range("M128").Value = 4750
$ws.Range("N128").Value = "`$/docena de matas"
$ws.Range("O128").Value = "Región Metropolitana"
$ws.Range("P128").Value = 792
$ws.Range("Q128").Value = 6
$ws.Range("D129").Value = (Get-Date -Year 2020 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I129").Value = "Segunda"
$ws.Range("J129").Value = 160
$ws.Range("K129").Value = 4000
$ws.Range("L129").Value = 4000
$ws.Range("M129").Value = 4000
$ws.Range("N129").Value = "`$/docena de matas"
$ws.Range("O129").Value = "Región Metropolitana"
$ws.Range("P129").Value = 667
$ws.Range("Q129").Value = 6
$ws.Range("D130").Value = (Get-Date -Year 2021 -Month 3 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I130").Value = "Primera"
$ws.Range("J130").Value = 340
$ws.Range("K130").Value = 3000
$ws.Range("L130").Value = 3500
$ws.Range("M130").Value = 3250
$ws.Range("N130").Value = "`$/docena de matas"
$ws.Range("O130").Value = "Región Metropolitana"
$ws.Range("P130").Value = 542
$ws.Range("Q130").Value = 6
$ws.Range("D131").Value = (Get-Date -Year 2020 -Month 11 -Day 24 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I131").Value = "Primera"
$ws.Range("J131").Value = 330
$ws.Range("K131").Value = 4500
$ws.Range("L131").Value = 5000
$ws.Range("M131").Value = 4773
$ws.Range("N131").Value = "`$/docena de matas"
$ws.Range("O131").Value = "Provincia de Chacabuco"
$ws.Range("P131").Value = 796
$ws.Range("Q131").Value = 6
$ws.Range("D132").Value = (Get-Date -Year 2021 -Month 5 -Day 17 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I132").Value = "Primera"
$ws.Range("J132").Value = 250
$ws.Range("K132").Value = 3000
$ws.Range("L132").Value = 3500
$ws.Range("M132").Value = 3250
$ws.Range("N132").Value = "`$/paquete"
$ws.Range("O132").Value = "Región de Arica y Parinacota"
$ws.Range("P132").Value = 3250
$ws.Range("Q132").Value = 1
$ws.Range("D133").Value = (Get-Date -Year 2021 -Month 1 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I133").Value = "Primera"
$ws.Range("J133").Value = 320
$ws.Range("K133").Value = 4000
$ws.Range("L133").Value = 4000
$ws.Range("M133").Value = 4000
$ws.Range("N133").Value = "`$/docena de matas"
$ws.Range("O133").Value = "Provincia de Chacabuco"
$ws.Range("P133").Value = 667
$ws.Range("Q133").Value = 6
$ws.Range("D134").Value = (Get-Date -Year 2020 -Month 12 -Day 23 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I134").Value = "Primera"
$ws.Range("J134").Value = 100
$ws.Range("K134").Value = 4000
$ws.Range("L134").Value = 4000
$ws.Range("M134").Value = 4000
$ws.Range("N134").Value = "`$/docena de matas"
$ws.Range("O134").Value = "Región Metropolitana"
$ws.Range("P134").Value = 667
$ws.Range("Q134").Value = 6
$ws.Range("D135").Value = (Get-Date -Year 2021 -Month 4 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I135").Value = "Primera"
$ws.Range("J135").Value = 250
$ws.Range("K135").Value = 5000
$ws.Range("L135").Value = 5000
$ws.Range("M135").Value = 5000
$ws.Range("N135").Value = "`$/docena de matas"
$ws.Range("O135").Value = "Región Metropolitana"
$ws.Range("P135").Value = 833
$ws.Range("Q135").Value = 6
$ws.Range("D136").Value = (Get-Date -Year 2021 -Month 4 -Day 30 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I136").Value = "Segunda"
$ws.Range("J136").Value = 160
$ws.Range("K136").Value = 4000
$ws.Range("L136").Value = 4000
$ws.Range("M136").Value = 4000
$ws.Range("N136").Value = "`$/docena de matas"
$ws.Range("O136").Value = "Región Metropolitana"
$ws.Range("P136").Value = 667
$ws.Range("Q136").Value = 6
$ws.Range("D137").Value = (Get-Date -Year 2021 -Month 2 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I137").Value = "Primera"
$ws.Range("J137").Value = 340
$ws.Range("K137").Value = 3000
$ws.Range("L137").Value = 3500
$ws.Range("M137").Value = 3250
$ws.Range("N137").Value = "`$/docena de matas"
$ws.Range("O137").Value = "Región Metropolitana"
$ws.Range("P137").Value = 542
$ws.Range("Q137").Value = 6
$ws.Range("D138").Value = (Get-Date -Year 2021 -Month 3 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I138").Value = "Primera"
$ws.Range("J138").Value = 340
$ws.Range("K138").Value = 3500
$ws.Range("L138").Value = 3500
$ws.Range("M138").Value = 3500
$ws.Range("N138").Value = "`$/docena de matas"
$ws.Range("O138").Value = "Región Metropolitana"
$ws.Range("P138").Value = 583
$ws.Range("Q138").Value = 6
$ws.Range("D139").Value = (Get-Date -Year 2021 -Month 3 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I139").Value = "Segunda"
$ws.Range("J139").Value = 160
$ws.Range("K139").Value = 2500
$ws.Range("L139").Value = 2500
$ws.Range("M139").Value = 2500
$ws.Range("N139").Value = "`$/docena de matas"
$ws.Range("O139").Value = "Región Metropolitana"
$ws.Range("P139").Value = 417
$ws.Range("Q139").Value = 6
$ws.Range("D140").Value = (Get-Date -Year 2020 -Month 11 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I140").Value = "Primera"
$ws.Range("J140").Value = 260
$ws.Range("K140").Value = 4000
$ws.Range("L140").Value = 5000
$ws.Range("M140").Value = 4615
$ws.Range("N140").Value = "`$/docena de matas"
$ws.Range("O140").Value = "Provincia de Chacabuco"
$ws.Range("P140").Value = 769
$ws.Range("Q140").Value = 6
$ws.Range("D141").Value = (Get-Date -Year 2020 -Month 12 -Day 21 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I141").Value = "Primera"
$ws.Range("J141").Value = 110
$ws.Range("K141").Value = 4000
$ws.Range("L141").Value = 4000
$ws.Range("M141").Value = 4000
$ws.Range("N141").Value = "`$/docena de matas"
$ws.Range("O141").Value = "Provincia de Chacabuco"
$ws.Range("P141").Value = 667
$ws.Range("Q141").Value = 6
$ws.Range("D142").Value = (Get-Date -Year 2021 -Month 1 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I142").Value = "Primera"
$ws.Range("J142").Value = 430
$ws.Range("K142").Value = 3000
$ws.Range("L142").Value = 4000
$ws.Range("M142").Value = 3500
$ws.Range("N142").Value = "`$/docena de matas"
$ws.Range("O142").Value = "Región Metropolitana"
$ws.Range("P142").Value = 583
$ws.Range("Q142").Value = 6
$ws.Range("D143").Value = (Get-Date -Year 2021 -Month 4 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I143").Value = "Primera"
$ws.Range("J143").Value = 340
$ws.Range("K143").Value = 3500
$ws.Range("L143").Value = 3500
$ws.Range("M143").Value = 3500
$ws.Range("N143").Value = "`$/docena de matas"
$ws.Range("O143").Value = "Región Metropolitana"
$ws.Range("P143").Value = 583
$ws.Range("Q143").Value = 6
$ws.Range("D144").Value = (Get-Date -Year 2021 -Month 4 -Day 15 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I144").Value = "Segunda"
$ws.Range("J144").Value = 160
$ws.Range("K144").Value = 2500
$ws.Range("L144").Value = 2500
$ws.Range("M144").Value = 2500
$ws.Range("N144").Value = "`$/docena de matas"
$ws.Range("O144").Value = "Región Metropolitana"
$ws.Range("P144").Value = 417
$ws.Range("Q144").Value = 6
$ws.Range("D145").Value = (Get-Date -Year 2021 -Month 3 -Day 17 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I145").Value = "Primera"
$ws.Range("J145").Value = 340
$ws.Range("K145").Value = 3500
$ws.Range("L145").Value = 3500
$ws.Range("M145").Value = 3500
$ws.Range("N145").Value = "`$/docena de matas"
$ws.Range("O145").Value = "Región Metropolitana"
$ws.Range("P145").Value = 583
$ws.Range("Q145").Value = 6
$ws.Range("D146").Value = (Get-Date -Year 2021 -Month 3 -Day 17 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I146").Value = "Segunda"
$ws.Range("J146").Value = 250
$ws.Range("K146").Value = 2500
$ws.Range("L146").Value = 2500
$ws.Range("M146").Value = 2500
$ws.Range("N146").Value = "`$/docena de matas"
$ws.Range("O146").Value = "Región Metropolitana"
$ws.Range("P146").Value = 417
$ws.Range("Q146").Value = 6
$ws.Range("D147").Value = (Get-Date -Year 2021 -Month 5 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I147").Value = "Primera"
$ws.Range("J147").Value = 250
$ws.Range("K147").Value = 3000
$ws.Range("L147").Value = 3000
$ws.Range("M147").Value = 3000
$ws.Range("N147").Value = "`$/paquete"
$ws.Range("O147").Value = "Región de Arica y Parinacota"
$ws.Range("P147").Value = 3000
$ws.Range("Q147").Value = 1
$ws.Range("D148").Value = (Get-Date -Year 2021 -Month 2 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I148").Value = "Primera"
$ws.Range("J148").Value = 390
$ws.Range("K148").Value = 3000
$ws.Range("L148").Value = 4000
$ws.Range("M148").Value = 3385
$ws.Range("N148").Value = "`$/docena de matas"
$ws.Range("O148").Value = "Provincia de Chacabuco"
$ws.Range("P148").Value = 564
$ws.Range("Q148").Value = 6
$ws.Range("D149").Value = (Get-Date -Year 2021 -Month 7 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I149").Value = "Primera"
$ws.Range("J149").Value = 250
$ws.Range("K149").Value = 4500
$ws.Range("L149").Value = 5000
$ws.Range("M149").Value = 4750
$ws.Range("N149").Value = "`$/paquete"
$ws.Range("O149").Value = "Región de Arica y Parinacota"
$ws.Range("P149").Value = 4750
$ws.Range("Q149").Value = 1
$ws.Range("D150").Value = (Get-Date -Year 2021 -Month 2 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I150").Value = "Primera"
$ws.Range("J150").Value = 580
$ws.Range("K150").Value = 2500
$ws.Range("L150").Value = 3000
$ws.Range("M150").Value = 2759
$ws.Range("N150").Value = "`$/docena de matas"
$ws.Range("O150").Value = "Provincia de Chacabuco"
$ws.Range("P150").Value = 460
$ws.Range("Q150").Value = 6
$ws.Range("D151").Value = (Get-Date -Year 2021 -Month 2 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I151").Value = "Primera"
$ws.Range("J151").Value = 150
$ws.Range("K151").Value = 3000
$ws.Range("L151").Value = 3000
$ws.Range("M151").Value = 3000
$ws.Range("N151").Value = "`$/docena de matas"
$ws.Range("O151").Value = "Región de O'Higgins"
$ws.Range("P151").Value = 500
$ws.Range("Q151").Value = 6
$ws.Range("D152").Value = (Get-Date -Year 2021 -Month 3 -Day 31 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I152").Value = "Primera"
$ws.Range("J152").Value = 340
$ws.Range("K152").Value = 3500
$ws.Range("L152").Value = 3500
$ws.Range("M152").Value = 3500
$ws.Range("N152").Value = "`$/docena de matas"
$ws.Range("O152").Value = "Región Metropolitana"
$ws.Range("P152").Value = 583
$ws.Range("Q152").Value = 6
$ws.Range("D153").Value = (Get-Date -Year 2021 -Month 3 -Day 31 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I153").Value = "Segunda"
$ws.Range("J153").Value = 160
$ws.Range("K153").Value = 2500
$ws.Range("L153").Value = 2500
$ws.Range("M153").Value = 2500
$ws.Range("N153").Value = "`$/docena de matas"
$ws.Range("O153").Value = "Región Metropolitana"
$ws.Range("P153").Value = 417
$ws.Range("Q153").Value = 6
$ws.Range("D154").Value = (Get-Date -Year 2021 -Month 1 -Day 21 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I154").Value = "Primera"
$ws.Range("J154").Value = 450
$ws.Range("K154").Value = 3000
$ws.Range("L154").Value = 4000
$ws.Range("M154").Value = 3556
$ws.Range("N154").Value = "`$/docena de matas"
$ws.Range("O154").Value = "Región Metropolitana"
$ws.Range("P154").Value = 593
$ws.Range("Q154").Value = 6
$ws.Range("D155").Value = (Get-Date -Year 2021 -Month 1 -Day 21 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I155").Value = "Primera"
$ws.Range("J155").Value = 150
$ws.Range("K155").Value = 3000
$ws.Range("L155").Value = 4000
$ws.Range("M155").Value = 3533
$ws.Range("N155").Value = "`$/docena de matas"
$ws.Range("O155").Value = "Región de O'Higgins"
$ws.Range("P155").Value = 589
$ws.Range("Q155").Value = 6
$ws.Range("D156").Value = (Get-Date -Year 2021 -Month 2 -Day 17 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I156").Value = "Primera"
$ws.Range("J156").Value = 520
$ws.Range("K156").Value = 3000
$ws.Range("L156").Value = 3000
$ws.Range("M156").Value = 3000
$ws.Range("N156").Value = "`$/docena de matas"
$ws.Range("O156").Value = "Región Metropolitana"
$ws.Range("P156").Value = 500
$ws.Range("Q156").Value = 6
$ws.Range("D157").Value = (Get-Date -Year 2021 -Month 5 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I157").Value = "Primera"
$ws.Range("J157").Value = 260
$ws.Range("K157").Value = 3000
$ws.Range("L157").Value = 3500
$ws.Range("M157").Value = 3250
$ws.Range("N157").Value = "`$/paquete"
$ws.Range("O157").Value = "Región de Arica y Parinacota"
$ws.Range("P157").Value = 3250
$ws.Range("Q157").Value = 1
$ws.Range("D158").Value = (Get-Date -Year 2021 -Month 1 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I158").Value = "Primera"
$ws.Range("J158").Value = 390
$ws.Range("K158").Value = 3000
$ws.Range("L158").Value = 4000
$ws.Range("M158").Value = 3821
$ws.Range("N158").Value = "`$/docena de matas"
$ws.Range("O158").Value = "Provincia de Chacabuco"
$ws.Range("P158").Value = 637
$ws.Range("Q158").Value = 6
$ws.Range("D159").Value = (Get-Date -Year 2021 -Month 1 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I159").Value = "Primera"
$ws.Range("J159").Value = 160
$ws.Range("K159").Value = 3000
$ws.Range("L159").Value = 4000
$ws.Range("M159").Value = 3500
$ws.Range("N159").Value = "`$/docena de matas"
$ws.Range("O159").Value = "Región de O'Higgins"
$ws.Range("P159").Value = 583
$ws.Range("Q159").Value = 6
$ws.Range("D160").Value = (Get-Date -Year 2021 -Month 6 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I160").Value = "Primera"
$ws.Range("J160").Value = 250
$ws.Range("K160").Value = 4500
$ws.Range("L160").Value = 5000
$ws.Range("M160").Value = 4750
$ws.Range("N160").Value = "`$/paquete"
$ws.Range("O160").Value = "Región de Arica y Parinacota"
$ws.Range("P160").Value = 4750
$ws.Range("Q160").Value = 1
$ws.Range("D161").Value = (Get-Date -Year 2021 -Month 1 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I161").Value = "Primera"
$ws.Range("J161").Value = 480
$ws.Range("K161").Value = 3000
$ws.Range("L161").Value = 4000
$ws.Range("M161").Value = 3417
$ws.Range("N161").Value = "`$/docena de matas"
$ws.Range("O161").Value = "Región Metropolitana"
$ws.Range("P161").Value = 570
$ws.Range("Q161").Value = 6
$ws.Range("D162").Value = (Get-Date -Year 2021 -Month 1 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I162").Value = "Primera"
$ws.Range("J162").Value = 110
$ws.Range("K162").Value = 3000
$ws.Range("L162").Value = 4000
$ws.Range("M162").Value = 3455
$ws.Range("N162").Value = "`$/docena de matas"
$ws.Range("O162").Value = "Región de O'Higgins"
$ws.Range("P162").Value = 576
$ws.Range("Q162").Value = 6
$ws.Range("D163").Value = (Get-Date -Year 2020 -Month 12 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I163").Value = "Primera"
$ws.Range("J163").Value = 60
$ws.Range("K163").Value = 5000
$ws.Range("L163").Value = 5000
$ws.Range("M163").Value = 5000
$ws.Range("N163").Value = "`$/docena de matas"
$ws.Range("O163").Value = "Provincia de Chacabuco"
$ws.Range("P163").Value = 833
$ws.Range("Q163").Value = 6
$ws.Range("D164").Value = (Get-Date -Year 2021 -Month 3 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I164").Value = "Primera"
$ws.Range("J164").Value = 250
$ws.Range("K164").Value = 3500
$ws.Range("L164").Value = 3500
$ws.Range("M164").Value = 3500
$ws.Range("N164").Value = "`$/docena de matas"
$ws.Range("O164").Value = "Región Metropolitana"
$ws.Range("P164").Value = 583
$ws.Range("Q164").Value = 6
$ws.Range("D165").Value = (Get-Date -Year 2021 -Month 3 -Day 29 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I165").Value = "Segunda"
$ws.Range("J165").Value = 160
$ws.Range("K165").Value = 2500
$ws.Range("L165").Value = 2500
$ws.Range("M165").Value = 2500
$ws.Range("N165").Value = "`$/docena de matas"
$ws.Range("O165").Value = "Región Metropolitana"
$ws.Range("P165").Value = 417
$ws.Range("Q165").Value = 6
$ws.Range("D166").Value = (Get-Date -Year 2020 -Month 12 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I166").Value = "Primera"
$ws.Range("J166").Value = 250
$ws.Range("K166").Value = 4500
$ws.Range("L166").Value = 5000
$ws.Range("M166").Value = 4700
$ws.Range("N166").Value = "`$/docena de matas"
$ws.Range("O166").Value = "Provincia de Chacabuco"
$ws.Range("P166").Value = 783
$ws.Range("Q166").Value = 6
$ws.Range("D167").Value = (Get-Date -Year 2020 -Month 12 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I167").Value = "Primera"
$ws.Range("J167").Value = 80
$ws.Range("K167").Value = 5000
$ws.Range("L167").Value = 5000
$ws.Range("M167").Value = 5000
$ws.Range("N167").Value = "`$/docena de matas"
$ws.Range("O167").Value = "Provincia de Quillota"
$ws.Range("P167").Value = 833
$ws.Range("Q167").Value = 6
